$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy number format/style from column W to the new V and X columns
$ws.Range("W2:W11").Copy()
$ws.Range("V2:V11").PasteSpecial(-4122)
$ws.Range("X2:X11").PasteSpecial(-4122)

# Header for new column X (shared string "11-jul")
$ws.Range("X1").Value = "11-jul"

# New data values for column V (between U and W)
$ws.Range("V2").Value = 13
$ws.Range("V3").Value = 24
$ws.Range("V4").Value = 9
$ws.Range("V5").Value = 10
$ws.Range("V6").Value = 13
$ws.Range("V7").Value = 15
$ws.Range("V8").Value = 13
$ws.Range("V9").Value = 19
$ws.Range("V10").Value = 20
$ws.Range("V11").Value = 10

# New data values for column X (after W)
$ws.Range("X2").Value = 12
$ws.Range("X3").Value = 18
$ws.Range("X4").Value = 7
$ws.Range("X5").Value = 9
$ws.Range("X6").Value = 13
$ws.Range("X7").Value = 15
$ws.Range("X8").Value = 15
$ws.Range("X9").Value = 10
$ws.Range("X10").Value = 21
$ws.Range("X11").Value = 22

[void]$ws.Range("X12").Select()

Write-Host "done"
